# Update countries & provincias Spain
# Applies the 30-Jul-2020 19:54 data refresh to the "Pais" sheet:
#  - updates the timestamp caption in A1
#  - refreshes the numeric COVID-19 stats (columns B-H) for the rows whose
#    underlying country ranking shifted
#  - swaps the country names for rows whose rank order changed so the row
#    with each country's updated stats also shows the right label

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp in title row
# Row 1: A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Julio de 2020 a las 19:54"

# Row 4: B4, C4, D4, E4, G4, H4
$ws.Range("B4").Value = 4600388
$ws.Range("C4").Value = 32351
$ws.Range("D4").Value = 2254135
$ws.Range("E4").Value = 2191712
$ws.Range("G4").Value = 701
$ws.Range("H4").Value = 154541

# Row 6: B6, C6, D6, E6, G6, H6
$ws.Range("B6").Value = 1639184
$ws.Range("C6").Value = 54800
$ws.Range("D6").Value = 1058659
$ws.Range("E6").Value = 544739
$ws.Range("G6").Value = 783
$ws.Range("H6").Value = 35786

# Row 20: B20, C20, D20, E20, G20, H20
$ws.Range("B20").Value = 229891
$ws.Range("C20").Value = 967
$ws.Range("D20").Value = 213539
$ws.Range("E20").Value = 10678
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 5674

# Row 21: B21, C21, E21, G21, H21
$ws.Range("B21").Value = 209501
$ws.Range("C21").Value = 690
$ws.Range("E21").Value = 8283
$ws.Range("G21").Value = 6
$ws.Range("H21").Value = 9218

# Row 36: A36, B36, C36, D36, E36, G36, H36
$ws.Range("A36").Value = "Israel"
$ws.Range("B36").Value = 69603
$ws.Range("C36").Value = 1304
$ws.Range("D36").Value = 35516
$ws.Range("E36").Value = 33588
$ws.Range("G36").Value = 8
$ws.Range("H36").Value = 499

# Row 37: A37, B37, C37, D37, E37, G37, H37
$ws.Range("A37").Value = "Ucrania"
$ws.Range("B37").Value = 68794
$ws.Range("C37").Value = 1197
$ws.Range("D37").Value = 38154
$ws.Range("E37").Value = 28967
$ws.Range("G37").Value = 23
$ws.Range("H37").Value = 1673

# Row 60: B60, C60, D60, E60, G60, H60
$ws.Range("B60").Value = 29831
$ws.Range("C60").Value = 602
$ws.Range("D60").Value = 20082
$ws.Range("E60").Value = 8549
$ws.Range("G60").Value = 14
$ws.Range("H60").Value = 1200

# Row 61: B61, C61, E61, G61, H61
$ws.Range("B61").Value = 26027
$ws.Range("C61").Value = 85
$ws.Range("E61").Value = 900
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 1763

# Row 64: A64, B64, C64, D64, E64, G64, H64
$ws.Range("A64").Value = "Marruecos"
$ws.Range("B64").Value = 23259
$ws.Range("C64").Value = 1046
$ws.Range("D64").Value = 17311
$ws.Range("E64").Value = 5602
$ws.Range("G64").Value = 12
$ws.Range("H64").Value = 346

# Row 65: A65, B65, C65, D65, E65, G65, H65
$ws.Range("A65").Value = "Uzbekistan"
$ws.Range("B65").Value = 23078
$ws.Range("C65").Value = 493
$ws.Range("D65").Value = 13432
$ws.Range("E65").Value = 9512
$ws.Range("G65").Value = 3
$ws.Range("H65").Value = 134

# Row 67: A67, B67, C67, D67, E67, G67, H67
$ws.Range("A67").Value = "Kenia"
$ws.Range("B67").Value = 19913
$ws.Range("C67").Value = 788
$ws.Range("D67").Value = 8121
$ws.Range("E67").Value = 11467
$ws.Range("G67").Value = 14
$ws.Range("H67").Value = 325

# Row 68: A68, B68, C68, D68, E68, G68, H68
$ws.Range("A68").Value = "Nepal"
$ws.Range("B68").Value = 19547
$ws.Range("C68").Value = 274
$ws.Range("D68").Value = 14248
$ws.Range("E68").Value = 5247
$ws.Range("G68").Value = 3
$ws.Range("H68").Value = 52

# Row 73: A73, B73, C73, D73, E73, G73, H73
$ws.Range("A73").Value = "Chequia"
$ws.Range("B73").Value = 16371
$ws.Range("C73").Value = 278
$ws.Range("D73").Value = 11482
$ws.Range("E73").Value = 4510
$ws.Range("G73").Value = 5
$ws.Range("H73").Value = 379

# Row 74: A74, B74, C74, D74, E74, G74, H74
$ws.Range("A74").Value = "Australia"
$ws.Range("B74").Value = 16303
$ws.Range("C74").Value = 721
$ws.Range("D74").Value = 10619
$ws.Range("E74").Value = 5495
$ws.Range("G74").Value = 13
$ws.Range("H74").Value = 189

# Row 75: A75, B75, C75, D75, E75, G75, H75
$ws.Range("A75").Value = "El Salvador"
$ws.Range("B75").Value = 16230
$ws.Range("C75").Value = 389
$ws.Range("D75").Value = 8206
$ws.Range("E75").Value = 7585
$ws.Range("G75").Value = 9
$ws.Range("H75").Value = 439

# Row 94: B94, C94, D94, E94
$ws.Range("B94").Value = 7242
$ws.Range("C94").Value = 59
$ws.Range("D94").Value = 6438
$ws.Range("E94").Value = 758

# Row 142: B142, C142, D142, E142, G142, H142
$ws.Range("B142").Value = 1181
$ws.Range("C142").Value = 2
$ws.Range("D142").Value = 667
$ws.Range("E142").Value = 441
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 73

# Row 184: A184, C184, E184, H184
$ws.Range("A184").Value = "Aruba"
$ws.Range("C184").Value = 1
$ws.Range("E184").Value = 12
$ws.Range("H184").Value = 3

# Row 185: A185, B185, D185, E185, H185
$ws.Range("A185").Value = "Monaco"
$ws.Range("B185").Value = 120
$ws.Range("D185").Value = 105
$ws.Range("E185").Value = 11
$ws.Range("H185").Value = 4

# Row 202: B202, C202, E202
$ws.Range("B202").Value = 25
$ws.Range("C202").Value = 1
$ws.Range("E202").Value = 3
